$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 0.7304773333333333
$ws.Range("H2").Value = 2.191432
$ws.Range("I2").Value = 0.03163269997405359
$ws.Range("J2").Value = 0.03163269997405359
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 3.752937333333333
$ws.Range("N2").Value = 11.258812
$ws.Range("O2").Value = 0.6855621274031838
$ws.Range("P2").Value = 0.6855621274031838
$ws.Range("Q2").Value = 2.741435655420444
$ws.Range("R2").Value = 24.672920898784
$ws.Range("S2").Value = 0.02168618108971882
$ws.Range("T2").Value = 0.02168618108971882
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 0.7304773333333333
$ws.Range("H3").Value = 2.191432
$ws.Range("I3").Value = 0.03163269997405359
$ws.Range("J3").Value = 0.03163269997405359
$ws.Range("O3").Value = 0.2368266084628361
$ws.Range("P3").Value = 0.2368266084628362
$ws.Range("Q3").Value = 0.9470256343528889
$ws.Range("R3").Value = 8.523230709176
$ws.Range("S3").Value = 0.007491465051377557
$ws.Range("T3").Value = 0.007491465051377559
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 0.7304773333333333
$ws.Range("H4").Value = 2.191432
$ws.Range("I4").Value = 0.03163269997405359
$ws.Range("J4").Value = 0.03163269997405359
$ws.Range("M4").Value = 0.4248633333333334
$ws.Range("N4").Value = 1.27459
$ws.Range("O4").Value = 0.07761126413398003
$ws.Range("P4").Value = 0.07761126413398005
$ws.Range("Q4").Value = 0.3103530347644445
$ws.Range("R4").Value = 2.79317731288
$ws.Range("S4").Value = 0.002455053832957217
$ws.Range("T4").Value = 0.002455053832957218
$ws.Range("I5").Value = 0.4074771110502447
$ws.Range("J5").Value = 0.4074771110502448
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 3.752937333333333
$ws.Range("N5").Value = 11.258812
$ws.Range("O5").Value = 0.6855621274031838
$ws.Range("P5").Value = 0.6855621274031838
$ws.Range("Q5").Value = 35.31384554328666
$ws.Range("R5").Value = 317.82460988958
$ws.Range("S5").Value = 0.2793508751197091
$ws.Range("T5").Value = 0.2793508751197092
$ws.Range("I6").Value = 0.4074771110502447
$ws.Range("J6").Value = 0.4074771110502448
$ws.Range("O6").Value = 0.2368266084628361
$ws.Range("P6").Value = 0.2368266084628362
$ws.Range("S6").Value = 0.0965014222362639
$ws.Range("T6").Value = 0.09650142223626393
$ws.Range("I7").Value = 0.4074771110502447
$ws.Range("J7").Value = 0.4074771110502448
$ws.Range("M7").Value = 0.4248633333333334
$ws.Range("N7").Value = 1.27459
$ws.Range("O7").Value = 0.07761126413398003
$ws.Range("P7").Value = 0.07761126413398005
$ws.Range("Q7").Value = 3.997817388816666
$ws.Range("R7").Value = 35.98035649935
$ws.Range("S7").Value = 0.03162481369427166
$ws.Range("T7").Value = 0.03162481369427167
$ws.Range("G8").Value = 12.95234266666667
$ws.Range("H8").Value = 38.857028
$ws.Range("I8").Value = 0.5608901889757016
$ws.Range("J8").Value = 0.5608901889757018
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 3.752937333333333
$ws.Range("N8").Value = 11.258812
$ws.Range("O8").Value = 0.6855621274031838
$ws.Range("P8").Value = 0.6855621274031838
$ws.Range("Q8").Value = 48.60933034785955
$ws.Range("R8").Value = 437.4839731307359
$ws.Range("S8").Value = 0.3845250711937558
$ws.Range("T8").Value = 0.3845250711937559
$ws.Range("G9").Value = 12.95234266666667
$ws.Range("H9").Value = 38.857028
$ws.Range("I9").Value = 0.5608901889757016
$ws.Range("J9").Value = 0.5608901889757018
$ws.Range("O9").Value = 0.2368266084628361
$ws.Range("P9").Value = 0.2368266084628362
$ws.Range("Q9").Value = 16.79203442806711
$ws.Range("R9").Value = 151.128309852604
$ws.Range("S9").Value = 0.1328337211751947
$ws.Range("T9").Value = 0.1328337211751947
$ws.Range("G10").Value = 12.95234266666667
$ws.Range("H10").Value = 38.857028
$ws.Range("I10").Value = 0.5608901889757016
$ws.Range("J10").Value = 0.5608901889757018
$ws.Range("M10").Value = 0.4248633333333334
$ws.Range("N10").Value = 1.27459
$ws.Range("O10").Value = 0.07761126413398003
$ws.Range("P10").Value = 0.07761126413398005
$ws.Range("Q10").Value = 5.502975479835556
$ws.Range("R10").Value = 49.52677931852001
$ws.Range("S10").Value = 0.04353139660675116
$ws.Range("T10").Value = 0.04353139660675117
